$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSubs = @("sub_006", "sub_007", "sub_008", "sub_009", "sub_010")
$row = 7
foreach ($sub in $newSubs) {
    $ws.Cells.Item($row, 1).Value = $sub
    $ws.Cells.Item($row, 2).Value = $false
    $row++
}

$ws.Range("B8").Select()
